# Weekly fruit/vegetable price update: insert 3 new daily records for
# "Vega Modelo de Temuco - Cebolla" dated 2021-11-09 (Excel serial 44509),
# pushing the existing rows 673:694 down to 676:697.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows right before the current row 673.
$ws.Rows("673:675").Insert()

# --- New row 673 ---
$ws.Range("A673").Value = 10
$ws.Range("B673").Value = "Vega Modelo de Temuco"
$ws.Range("C673").Value = "La Araucanía"
$ws.Range("D673").Value = 44509
$ws.Range("E673").Value = 9
$ws.Range("F673").Value = 100112004
$ws.Range("G673").Value = "Cebolla"
$ws.Range("H673").Value = "Morada(o)"
$ws.Range("I673").Value = "Primera"
$ws.Range("J673").Value = 50
$ws.Range("K673").Value = 8000
$ws.Range("L673").Value = 8000
$ws.Range("M673").Value = 8000
$ws.Range("N673").Value = "$/malla 18 kilos"
$ws.Range("O673").Value = "Perú"
$ws.Range("P673").Value = 444
$ws.Range("Q673").Value = 18
$ws.Range("R673").Value = "Hortaliza"

# --- New row 674 ---
$ws.Range("A674").Value = 10
$ws.Range("B674").Value = "Vega Modelo de Temuco"
$ws.Range("C674").Value = "La Araucanía"
$ws.Range("D674").Value = 44509
$ws.Range("E674").Value = 9
$ws.Range("F674").Value = 100112004
$ws.Range("G674").Value = "Cebolla"
$ws.Range("H674").Value = "Sin especificar"
$ws.Range("I674").Value = "1a nueva(o)"
$ws.Range("J674").Value = 300
$ws.Range("K674").Value = 5000
$ws.Range("L674").Value = 5000
$ws.Range("M674").Value = 5000
$ws.Range("N674").Value = "$/malla 18 kilos"
$ws.Range("O674").Value = "Región Metropolitana"
$ws.Range("P674").Value = 278
$ws.Range("Q674").Value = 18
$ws.Range("R674").Value = "Hortaliza"

# --- New row 675 ---
$ws.Range("A675").Value = 10
$ws.Range("B675").Value = "Vega Modelo de Temuco"
$ws.Range("C675").Value = "La Araucanía"
$ws.Range("D675").Value = 44509
$ws.Range("E675").Value = 9
$ws.Range("F675").Value = 100112004
$ws.Range("G675").Value = "Cebolla"
$ws.Range("H675").Value = "Sin especificar"
$ws.Range("I675").Value = "Primera"
$ws.Range("J675").Value = 1000
$ws.Range("K675").Value = 4500
$ws.Range("L675").Value = 4500
$ws.Range("M675").Value = 4500
$ws.Range("N675").Value = "$/malla 18 kilos"
$ws.Range("O675").Value = "Perú"
$ws.Range("P675").Value = 250
$ws.Range("Q675").Value = 18
$ws.Range("R675").Value = "Hortaliza"
